$d = $word.ActiveDocument

$replacements = @(
    @("286÷4=", "125÷4="),
    @("828÷5=", "854÷6="),
    @("424÷4=", "542÷5="),
    @("444÷8=", "755÷7="),
    @("263÷6=", "566÷5="),
    @("230÷7=", "310÷7="),
    @("171÷9=", "701÷5="),
    @("559÷6=", "794÷9="),
    @("725÷3=", "299÷2="),
    @("291÷3=", "826÷2="),
    @("459÷7=", "126÷5="),
    @("871÷8=", "952÷9="),
    @("813÷4=", "396÷7="),
    @("905÷9=", "390÷2="),
    @("981÷2=", "340÷3="),
    @("939÷3=", "925÷3="),
    @("490÷3=", "157÷2="),
    @("589÷4=", "396÷3="),
    @("370÷3=", "559÷7="),
    @("162÷6=", "734÷6="),
    @("351÷7=", "237÷2="),
    @("666÷5=", "469÷2="),
    @("829÷8=", "483÷4="),
    @("573÷2=", "912÷6="),
    @("323÷4=", "223÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
